$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 497.7
$ws.Range("I33").Value = 580.1429000000001
$ws.Range("J33").Value = 305.33334
$ws.Range("K33").Value = 580.1429000000001
$ws.Range("L33").Value = 305.33334
$ws.Range("M33").Value = -351.1429000000001
$ws.Range("N33").Value = -763.33334
$ws.Range("H62").Value = 6682.143
$ws.Range("I62").Value = 6682.143
$ws.Range("K62").Value = 6682.143
$ws.Range("M62").Value = -6058.143
$ws.Range("H65").Value = 6682.143
$ws.Range("I65").Value = 6682.143
$ws.Range("K65").Value = 33410.715
$ws.Range("M65").Value = -30290.715

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2357889.8
$ws.Range("I2").Value = 3536076.5
$ws.Range("K2").Value = 3536076.5
$ws.Range("M2").Value = -3535963.5
$ws.Range("H32").Value = 12861.3
$ws.Range("I32").Value = 8147.5
$ws.Range("K32").Value = 8147.5
$ws.Range("M32").Value = -7860.5
$ws.Range("H63").Value = 4969.1113
$ws.Range("I63").Value = 2455.4
$ws.Range("K63").Value = 2455.4
$ws.Range("M63").Value = -1769.4
$ws.Range("H66").Value = 4969.1113
$ws.Range("I66").Value = 2455.4
$ws.Range("K66").Value = 12277
$ws.Range("M66").Value = -8845
$ws.Range("H92").Value = 18000
$ws.Range("J92").Value = 18000
$ws.Range("L92").Value = 18000
$ws.Range("N92").Value = -22992
$ws.Range("H110").Value = 713319
$ws.Range("I110").Value = 751833.6
$ws.Range("K110").Value = 751833.6
$ws.Range("M110").Value = -749788.6
$ws.Range("H116").Value = 2357889.8
$ws.Range("I116").Value = 3536076.5
$ws.Range("K116").Value = 3536076.5
$ws.Range("M116").Value = -3533782.5
$ws.Range("H122").Value = 466385.03
$ws.Range("I122").Value = 2201.4482
$ws.Range("K122").Value = 6604.344599999999
$ws.Range("M122").Value = -4154.344599999999
$ws.Range("H132").Value = 2211.3257
$ws.Range("I132").Value = 1751.825
$ws.Range("K132").Value = 5255.475
$ws.Range("M132").Value = -2725.475

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2357889.8
$ws.Range("I3").Value = 3536076.5
$ws.Range("K3").Value = 3536076.5
$ws.Range("M3").Value = -3535962.5
$ws.Range("H80").Value = 456.32352
$ws.Range("I80").Value = 496
$ws.Range("J80").Value = 434.68182
$ws.Range("K80").Value = 496
$ws.Range("L80").Value = 434.68182
$ws.Range("M80").Value = 502
$ws.Range("N80").Value = -2430.68182
$ws.Range("H83").Value = 456.32352
$ws.Range("I83").Value = 496
$ws.Range("J83").Value = 434.68182
$ws.Range("K83").Value = 2480
$ws.Range("L83").Value = 2173.4091
$ws.Range("M83").Value = 2512
$ws.Range("N83").Value = -12157.4091
$ws.Range("H127").Value = 45666.668
$ws.Range("J127").Value = 45666.668
$ws.Range("L127").Value = 45666.668
$ws.Range("N127").Value = -55586.668

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25790.023
$ws.Range("I31").Value = 3333.8572
$ws.Range("K31").Value = 3333.8572
$ws.Range("M31").Value = -3038.8572
$ws.Range("H34").Value = 25790.023
$ws.Range("I34").Value = 3333.8572
$ws.Range("K34").Value = 3333.8572
$ws.Range("M34").Value = -3131.8572
$ws.Range("H58").Value = 5714.1333
$ws.Range("I58").Value = 6452.6665
$ws.Range("J58").Value = 3990.889
$ws.Range("K58").Value = 6452.6665
$ws.Range("L58").Value = 3990.889
$ws.Range("M58").Value = -6249.6665
$ws.Range("N58").Value = -4396.889
$ws.Range("H122").Value = 2295.7778
$ws.Range("I122").Value = 2136.7058
$ws.Range("K122").Value = 6410.117400000001
$ws.Range("M122").Value = -3960.117400000001
$ws.Range("H132").Value = 55987.445
$ws.Range("I132").Value = 38630.855
$ws.Range("K132").Value = 115892.565
$ws.Range("M132").Value = -113362.565
$ws.Range("H134").Value = 3401.1875
$ws.Range("I134").Value = 2306.762
$ws.Range("K134").Value = 6920.286
$ws.Range("M134").Value = -4385.286
$ws.Range("H136").Value = 5714.1333
$ws.Range("I136").Value = 6452.6665
$ws.Range("J136").Value = 3990.889
$ws.Range("K136").Value = 19357.9995
$ws.Range("L136").Value = 11972.667
$ws.Range("M136").Value = -16807.9995
$ws.Range("N136").Value = -17072.667

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 306.66666
$ws.Range("I17").Value = 133.33333
$ws.Range("J17").Value = 393.33334
$ws.Range("K17").Value = 399.99999
$ws.Range("L17").Value = 1180.00002
$ws.Range("M17").Value = -230.99999
$ws.Range("N17").Value = -1518.00002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 25002982
$ws.Range("I2").Value = 550.04346
$ws.Range("K2").Value = 550.04346
$ws.Range("M2").Value = -437.04346
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H80").Value = 47683196
$ws.Range("I80").Value = 74928760
$ws.Range("J80").Value = 3450
$ws.Range("K80").Value = 74928760
$ws.Range("L80").Value = 3450
$ws.Range("M80").Value = -74927762
$ws.Range("N80").Value = -5446
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H83").Value = 47683196
$ws.Range("I83").Value = 74928760
$ws.Range("J83").Value = 3450
$ws.Range("K83").Value = 374643800
$ws.Range("L83").Value = 17250
$ws.Range("M83").Value = -374638808
$ws.Range("N83").Value = -27234
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H126").Value = 3792177
$ws.Range("I126").Value = 7578582.5
$ws.Range("K126").Value = 22735747.5
$ws.Range("M126").Value = -22733277.5
$ws.Range("H132").Value = 2980.318
$ws.Range("I132").Value = 2731.25
$ws.Range("K132").Value = 8193.75
$ws.Range("M132").Value = -5663.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8874.777
$ws.Range("I7").Value = 6375.25
$ws.Range("K7").Value = 6375.25
$ws.Range("M7").Value = -6263.25
$ws.Range("H46").Value = 7081.091
$ws.Range("J46").Value = 8749
$ws.Range("L46").Value = 8749
$ws.Range("N46").Value = -9125
$ws.Range("H126").Value = 8874.777
$ws.Range("I126").Value = 6375.25
$ws.Range("K126").Value = 19125.75
$ws.Range("M126").Value = -16655.75
$ws.Range("H132").Value = 16936.223
$ws.Range("I132").Value = 17803.25
$ws.Range("K132").Value = 53409.75
$ws.Range("M132").Value = -50879.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 9589
$ws.Range("I55").Value = 4000
$ws.Range("J55").Value = 10986.25
$ws.Range("K55").Value = 4000
$ws.Range("L55").Value = 10986.25
$ws.Range("M55").Value = -3723
$ws.Range("N55").Value = -11540.25
$ws.Range("H122").Value = 6719.6
$ws.Range("I122").Value = 4534.3335
$ws.Range("K122").Value = 13603.0005
$ws.Range("M122").Value = -11153.0005
$ws.Range("H132").Value = 12639097
$ws.Range("I132").Value = 15387937
$ws.Range("K132").Value = 46163811
$ws.Range("M132").Value = -46161281
$ws.Range("H136").Value = 3716.7188
$ws.Range("I136").Value = 3391.9583
$ws.Range("J136").Value = 4691
$ws.Range("K136").Value = 10175.8749
$ws.Range("L136").Value = 14073
$ws.Range("M136").Value = -7625.874899999999
$ws.Range("N136").Value = -19173
